$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.598.40'
$ws.Range("E2").Value = '  -3.37%  '
$ws.Range("D3").Value = '3.335.66'
$ws.Range("E3").Value = '  -3.87%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '548.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.611'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.12%  '
$ws.Range("D8").Value = '3.326.97'
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  -3.01%  '
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.67%  '
$ws.Range("D15").Value = '3.869.91'
$ws.Range("E15").Value = '  -4.03%  '
$ws.Range("D16").Value = '3.348.60'
$ws.Range("E16").Value = '  -3.65%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '17.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.57%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.117'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.08%  '
$ws.Range("D20").Value = '63.600.28'
$ws.Range("E20").Value = '  -3.46%  '
$ws.Range("E21").Value = '  -1.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '415.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("E24").Value = '  +4.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '83.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.51%  '
$ws.Range("E28").Value = '  -5.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.58'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.49%  '
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '574.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.07%  '
$ws.Range("E34").Value = '  -3.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.60'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.92%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.147'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.43%  '
$ws.Range("D39").Value = '0.0₃0738'
$ws.Range("E39").Value = '  -6.34%  '
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("E41").Value = '  -3.94%  '
$ws.Range("D42").Value = '3.142.99'
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("E47").Value = '  -5.91%  '
$ws.Range("E48").Value = '  -4.62%  '
$ws.Range("E49").Value = '  -3.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.34%  '
$ws.Range("E51").Value = '  -5.03%  '
